# Adds the new milk-product order lines (rows 4-7) to the active sheet,
# matching the "Hillcrest Dairy" order layout: SKU | Name | Quantity | Cost Per | Total Cost.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("wholeg", "Milk - Whole",               "12", "16.08", "192.96"),
    @("skimg",  "Milk - Skim",                "1",  "15.56", "15.56"),
    @("twog",   "Milk - 2%",                  "20", "16.08", "321.60"),
    @("choqt",  "Milk - Chocolate (9/32oz)",  "1",  "16.65", "16.65")
)

$startRow = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $item = $rows[$i]

    # SKU / Name are plain text already.
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]

    # Quantity / Cost Per / Total Cost look numeric, but the sheet stores
    # them as text (matching the existing rows) — force text format first
    # so Excel doesn't auto-coerce the values to numbers, then clear the
    # formatting back off so the cell ends up plain text with no explicit
    # number-format style applied (matching the rest of the sheet).
    $qtyCell = $ws.Cells.Item($r, 3)
    $qtyCell.NumberFormat = "@"
    $qtyCell.Value = $item[2]
    $qtyCell.ClearFormats()

    $costCell = $ws.Cells.Item($r, 4)
    $costCell.NumberFormat = "@"
    $costCell.Value = $item[3]
    $costCell.ClearFormats()

    $totalCell = $ws.Cells.Item($r, 5)
    $totalCell.NumberFormat = "@"
    $totalCell.Value = $item[4]
    $totalCell.ClearFormats()
}
